$d = $word.ActiveDocument

# Locate the end of the existing "service call" message (end of the
# first invalid-query run) using Find, which yields a properly anchored
# Range we can reliably use for insertion.
$searchRange = $d.Content
$searchRange.Find.Execute("service call", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$text1 = "    <---Invalid query statement: Couldn't find the 'self' variable"
$text2 = "    <---Invalid query statement: missing feature access or service call"

# Insert first new message right after the found text, then format it.
$insertStart1 = $searchRange.End
$searchRange.InsertAfter($text1)
$newRange1 = $d.Range($insertStart1, $insertStart1 + $text1.Length)
$newRange1.Font.Bold = $true
$newRange1.Font.Color = 255

# Insert the second new message right after the first one, then format it.
$insertStart2 = $newRange1.End
$newRange1.InsertAfter($text2)
$newRange2 = $d.Range($insertStart2, $insertStart2 + $text2.Length)
$newRange2.Font.Bold = $true
$newRange2.Font.Color = 255
